# Applies the "moved more features from the bug tracker" edit:
#  1. Update the DATE field result text (Tuesday, 14 April 2009 -> Wednesday, 15 April 2009)
#  2. "What-if" -> "What-if / Concept exploration" (as a separate trailing run)
#  3. "Goals" -> "Goals:" (as a separate trailing run)
#  4. New "Experiments:" section with two bullet points inserted after the Goals bullet
#  5. "Parameters" -> "Parameters:" (as a separate trailing run) and a lastRenderedPageBreak
#     marker added before the first bullet under it
#  6. "File I/O" -> "File I/O:" (as a separate trailing run)

$d = $word.ActiveDocument
$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-PackageXml([string]$bodyXml) {
    return '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
           $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Get-ParaByText([string]$text) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $text) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1. Date field
# ---------------------------------------------------------------------------
foreach ($f in $d.Fields) {
    $f.Result.Text = "Wednesday, 15 April 2009"
}

# ---------------------------------------------------------------------------
# 2. "What-if" -> "What-if" + " / Concept exploration" (two runs)
# ---------------------------------------------------------------------------
$p = Get-ParaByText "What-if"
$xml = New-PackageXml(
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:r><w:t>What-if</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> / Concept exploration</w:t></w:r></w:p>'
)
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 3. "Goals" -> "Goals" + ":" (two runs)
# ---------------------------------------------------------------------------
$p = Get-ParaByText "Goals"
$xml = New-PackageXml(
    '<w:p><w:r><w:t>Goals</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>'
)
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 4. New "Experiments:" block inserted after the Goals bullet + blank paragraph
# ---------------------------------------------------------------------------
# Note: the immediately-following (pre-existing) blank paragraph is left
# untouched and serves as the trailing blank line before "Parameters" - if an
# empty <w:p/> is put at the end of this replacement body, it swallows
# whatever paragraph comes right after it, so we deliberately stop the body
# at the non-blank "boundaries of the geometry." paragraph.
$p = Get-ParaByText "Allow users to specify the desired error bounds on the goal."
$nextP = $p.Next()
$body = '<w:p/>' +
    '<w:p><w:r><w:t>Experiments:</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:r><w:t>The user should only have to specify the physical properties of the experiment.  The system should then be able to suggest which models / boundary conditions are applicable. Obviously the user should be able to amend these suggestions.</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Users should only have to specify the geometry of the problem, not the additional boundaries. The only problem will lie in the differentiation of internal / </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>external  &amp;</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> boundaries of the geometry.</w:t></w:r></w:p>'
$xml = New-PackageXml($body)
$nextP.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 5. "Parameters" -> "Parameters" + ":" and add a lastRenderedPageBreak marker
#    before the first bullet beneath it
# ---------------------------------------------------------------------------
$p = Get-ParaByText "Parameters"
$xml = New-PackageXml(
    '<w:p><w:r><w:t>Parameters</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>'
)
$p.Range.InsertXML($xml)

$p = Get-ParaByText "Parameters should always have units"
$xml = New-PackageXml(
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:r><w:lastRenderedPageBreak/><w:t>Parameters should always have units</w:t></w:r></w:p>'
)
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 6. "File I/O" -> "File I/O" + ":" (two runs)
# ---------------------------------------------------------------------------
$p = Get-ParaByText "File I/O"
$xml = New-PackageXml(
    '<w:p><w:r><w:t>File I/O</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>'
)
$p.Range.InsertXML($xml)
